# Add a new row (row 25) to the SQL50 pattern tracker table on Sheet1,
# documenting LeetCode 1045 "Customers Who Bought All Products".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$question    = "1045. Customers Who Bought All Products"
$difficulty  = "Medium"
$pattern     = "Sorting and Grouping"
$notes       = "select having count(distinct product_key) = (select count(*) from Product)"
$link        = "https://leetcode.com/problems/customers-who-bought-all-products/solutions/3865452/100-easy-fast-clean-solution/?envType=study-plan-v2&envId=top-sql-50 "

$ws.Range("A25").Value = $question
$ws.Range("B25").Value = $difficulty
$ws.Range("C25").Value = $pattern
$ws.Range("D25").Value = $notes
$ws.Range("E25").Value = $link

# Match the "Medium" fill styling already used on the other Medium rows (e.g. B16)
$ws.Range("B25").Interior.Color = $ws.Range("B16").Interior.Color

# Add the hyperlink on the new Link cell, matching the style of existing rows
$ws.Hyperlinks.Add($ws.Range("E25"), $link) | Out-Null
$ws.Range("E25").Style = $ws.Range("E24").Style

# Expand the table (ListObject) range to include the new row
$table = $ws.ListObjects.Item(1)
$table.Resize($ws.Range("A1:E25"))

$ws.Range("E32").Select() | Out-Null
